# Generate Report for Handoff
#
# A new handoff xliff was generated, so the GUID-named source/handoff file
# references (and their generation timestamps) are refreshed across the
# three worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$oldGuid = "9a07edeb-751d-4068-983a-5d2e2ed23219"
$newGuid = "cdd2fe4f-bb6d-42c0-b3e7-1578ca84ad4c"

$newMd = "$newGuid.md"
$newMdPath = "e2e\$newGuid.md"

$newZhXlf = "$newGuid.0876aadc5196a6a07f6a2645304a0c90a5b2653b.zh-cn.xlf"
$newDeXlf = "$newGuid.0876aadc5196a6a07f6a2645304a0c90a5b2653b.de-de.xlf"

# The external hyperlink target is unchanged by this edit - only the
# friendly display text (which mirrors the file name) is refreshed.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4625a941de1a5551d4d4a9b69e8ca29fd79c3aa4/e2e/$oldGuid.md"

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkAddress, [Type]::Missing, [Type]::Missing, $newMdPath)
$wsOverview.Range("G2").Value = "2016-11-03 19:45:52"

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $linkAddress, [Type]::Missing, [Type]::Missing, $newMd)
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = "2016-11-03 19:45:39"

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $linkAddress, [Type]::Missing, [Type]::Missing, $newMd)
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = "2016-11-03 19:45:52"
